$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "testx.pdf"
$ws.Range("A7").Value = "test.docx"

$ws.Range("A7").Select()
